$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Homeless Inc."
$ws.Range("B3").Value = "Дубов Александр Лесович"
$ws.Range("B4").Value = "Директор"
$ws.Range("B6").Value = "dubov@mail.ru"
$ws.Range("E8").Value = "Pushkina 88/4"
$ws.Range("E9").Value = "Moscow"
$ws.Range("E10").Value = "www.MMM.ru"
